# Refreshed crypto-market snapshot: updates the Price (D) and 1h Volume
# change (E) columns in place. Column D stores plain text in the source
# file (t="inlineStr") even when a price looks like a bare number, so for
# those values we lead with a text-marker apostrophe to stop Excel from
# auto-converting the cell to a Number (same as typing '167.81 by hand).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.925.05'
$ws.Range("E2").Value = '  +2.28%  '
$ws.Range("D3").Value = '3.742.54'
$ws.Range("E3").Value = '  +0.61%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '''601.53'
$ws.Range("E5").Value = '  +2.02%  '
$ws.Range("D6").Value = '''167.81'
$ws.Range("E6").Value = '  -2.02%  '
$ws.Range("D7").Value = '3.742.63'
$ws.Range("E7").Value = '  +0.58%  '
$ws.Range("D9").Value = '''0.534'
$ws.Range("E9").Value = '  +3.49%  '
$ws.Range("E10").Value = '  +5.16%  '
$ws.Range("E11").Value = '  +2.98%  '
$ws.Range("E12").Value = '  +0.40%  '
$ws.Range("D13").Value = '''38.19'
$ws.Range("E13").Value = '  +2.37%  '
$ws.Range("E14").Value = '  +1.86%  '
$ws.Range("D15").Value = '4.368.66'
$ws.Range("E15").Value = '  +0.73%  '
$ws.Range("D16").Value = '3.746.53'
$ws.Range("E16").Value = '  +0.98%  '
$ws.Range("D17").Value = '68.918.00'
$ws.Range("E17").Value = '  +2.35%  '
$ws.Range("D18").Value = '''7.26'
$ws.Range("E18").Value = '  +2.23%  '
$ws.Range("E19").Value = '  +1.10%  '
$ws.Range("D20").Value = '''17.24'
$ws.Range("E20").Value = '  +7.78%  '
$ws.Range("D21").Value = '''497.00'
$ws.Range("E21").Value = '  +2.79%  '
$ws.Range("E22").Value = '  +14.84%  '
$ws.Range("E23").Value = '  +2.24%  '
$ws.Range("D24").Value = '''85.23'
$ws.Range("E24").Value = '  +2.44%  '
$ws.Range("E25").Value = '  -0.86%  '
$ws.Range("E26").Value = '  +2.86%  '
$ws.Range("D27").Value = '''12.27'
$ws.Range("E27").Value = '  +1.50%  '
$ws.Range("D28").Value = '''10.13'
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("E29").Value = '  +0.25%  '
$ws.Range("E30").Value = '  +1.99%  '
$ws.Range("E31").Value = '  +2.18%  '
$ws.Range("D32").Value = '''7.92'
$ws.Range("E32").Value = '  +3.33%  '
$ws.Range("D33").Value = '''31.75'
$ws.Range("E33").Value = '  -1.01%  '
$ws.Range("D34").Value = '3.887.37'
$ws.Range("E34").Value = '  +1.00%  '
$ws.Range("E35").Value = '  +1.42%  '
$ws.Range("D36").Value = '3.675.59'
$ws.Range("E36").Value = '  +0.53%  '
$ws.Range("E37").Value = '  +0.10%  '
$ws.Range("E38").Value = '  +2.43%  '
$ws.Range("E39").Value = '  +2.54%  '
$ws.Range("E40").Value = '  +0.43%  '
$ws.Range("D41").Value = '''0.324'
$ws.Range("E41").Value = '  +1.22%  '
$ws.Range("D42").Value = '''438.41'
$ws.Range("E42").Value = '  -1.48%  '
$ws.Range("D43").Value = '''49.02'
$ws.Range("E43").Value = '  +0.59%  '
$ws.Range("E44").Value = '  +1.02%  '
$ws.Range("E45").Value = '  +2.07%  '
$ws.Range("E46").Value = '  +2.56%  '
$ws.Range("D48").Value = '''40.43'
$ws.Range("E48").Value = '  -1.54%  '
$ws.Range("D49").Value = '''142.61'
$ws.Range("E49").Value = '  +1.31%  '
$ws.Range("E50").Value = '  +2.43%  '
$ws.Range("D51").Value = '2.748.79'
$ws.Range("E51").Value = '  -1.07%  '
